{"js": "// Replace the closing sentence of the final paragraph with the new,\n// expanded explanation about the count functions / constraints results.\nconst oldTail = \"Unfortunately we didn\\u2019t find the time to back up this claim.\";\nconst newTail =\n  \"Though not a very strong test case, we made count functions which counts the amount of constraints on both hard and easy Sudokus. \" +\n  \"The results gave that hard Sudokus have on average about 20 constraints more. This means that in the start situation the hard Sudokus \" +\n  \"have on 59 free positions there are on average 20 more possibilities to enter in these positions. More posibilities normally means that more thinking is needed to solve this. \";\n\nconst results = context.document.body.search(oldTail, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \" + oldTail);\n}\n\nresults.items[0].insertText(newTail, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the closing sentence of the final paragraph with the new,\n# expanded explanation about the count functions / constraints results.\n$d = $word.ActiveDocument\n\n$oldTail = \"Unfortunately we didn\" + [char]0x2019 + \"t find the time to back up this claim.\"\n$newTail = \"Though not a very strong test case, we made count functions which counts the amount of constraints on both hard and easy Sudokus. \" + `\n    \"The results gave that hard Sudokus have on average about 20 constraints more. This means that in the start situation the hard Sudokus \" + `\n    \"have on 59 free positions there are on average 20 more possibilities to enter in these positions. More posibilities normally means that more thinking is needed to solve this. \"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldTail\n$find.Replacement.Text = $newTail\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
